# Apply corrected Morph values (column I) for buzzard rows, per commit
# "inserted correct Morphs of buzzards".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 had no confirmed morph reading -> clear the previously-guessed value.
$ws.Range("I51").ClearContents()

# Corrected morph codes for the remaining rows (changed from the placeholder 3).
$morphUpdates = @{
    41 = 4
    42 = 4
    44 = 2
    45 = 2
    46 = 4
    47 = 4
    48 = 2
    50 = 2
    52 = 4
    55 = 2
    56 = 2
    57 = 2
    61 = 2
    62 = 4
    64 = 2
    67 = 2
    71 = 2
    72 = 4
    73 = 4
    75 = 2
    76 = 2
    78 = 2
    79 = 2
    80 = 2
    81 = 2
    82 = 4
    85 = 2
    86 = 2
    87 = 4
    88 = 4
    89 = 2
}

foreach ($row in $morphUpdates.Keys) {
    $ws.Range("I$row").Value = $morphUpdates[$row]
}

# Move the view/selection to where the last edit (row 51) was made.
$ws.Range("I51").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
